$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18 (G18 old context)
$ws.Range("H18").Value = 1245.5714
$ws.Range("I18").Value = 469.83334
$ws.Range("K18").Value = 469.83334
$ws.Range("M18").Value = -185.83334

# Row 135 (G135 old context)
$ws.Range("H135").Value = 17858380
$ws.Range("I135").Value = 1230.6
$ws.Range("J135").Value = 62501250
$ws.Range("K135").Value = 11075.4
$ws.Range("L135").Value = 562511250
$ws.Range("M135").Value = -8540.4
$ws.Range("N135").Value = -562516320

# Row 137 (G137 old context)
$ws.Range("H137").Value = 4011.2666
$ws.Range("I137").Value = 1368.4445
$ws.Range("K137").Value = 4105.333500000001
$ws.Range("M137").Value = -1555.333500000001

# Row 141 (G141 old context)
$ws.Range("H141").Value = 10708.272
$ws.Range("I141").Value = 13460.875
$ws.Range("J141").Value = 3368
$ws.Range("K141").Value = 40382.625
$ws.Range("L141").Value = 10104
$ws.Range("M141").Value = -35202.625
$ws.Range("N141").Value = -20464

$ws = $wb.Worksheets.Item("ARM")
# Row 39 (G39 old context)
$ws.Range("H39").Value = 19749.25
$ws.Range("I39").Value = 19000
$ws.Range("J39").Value = 19999
$ws.Range("K39").Value = 19000
$ws.Range("L39").Value = 19999
$ws.Range("M39").Value = -18480
$ws.Range("N39").Value = -21039

# Row 45 (G45 old context)
$ws.Range("H45").Value = 4018.7273
$ws.Range("I45").Value = 3304.25
$ws.Range("J45").Value = 4427
$ws.Range("K45").Value = 3304.25
$ws.Range("L45").Value = 4427
$ws.Range("M45").Value = -2927.25
$ws.Range("N45").Value = -5181

# Row 49 (G49 old context)
$ws.Range("H49").Value = 14076.923
$ws.Range("J49").Value = 14076.923
$ws.Range("L49").Value = 14076.923
$ws.Range("N49").Value = -14596.923

# Row 61 (G61 old context)
$ws.Range("H61").Value = 22729380
$ws.Range("I61").Value = 29413746
$ws.Range("K61").Value = 29413746
$ws.Range("M61").Value = -29413534

# Row 63 (G63 old context)
$ws.Range("H63").Value = 9600
$ws.Range("I63").Value = 8000
$ws.Range("K63").Value = 8000
$ws.Range("M63").Value = -7314

# Row 66 (G66 old context)
$ws.Range("H66").Value = 9600
$ws.Range("I66").Value = 8000
$ws.Range("K66").Value = 40000
$ws.Range("M66").Value = -36568

# Row 132 (G132 old context)
$ws.Range("H132").Value = 21308350
$ws.Range("I132").Value = 1967.1613
$ws.Range("J132").Value = 62589468
$ws.Range("K132").Value = 5901.4839
$ws.Range("L132").Value = 187768404
$ws.Range("M132").Value = -3371.4839
$ws.Range("N132").Value = -187773464

# Row 136 (G136 old context)
$ws.Range("H136").Value = 22729380
$ws.Range("I136").Value = 29413746
$ws.Range("K136").Value = 88241238
$ws.Range("M136").Value = -88238688

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (G20 old context)
$ws.Range("H20").Value = 4312.375
$ws.Range("I20").Value = 3495
$ws.Range("K20").Value = 3495
$ws.Range("M20").Value = -3248

# Row 134 (G134 old context)
$ws.Range("H134").Value = 2511.0293
$ws.Range("I134").Value = 2463.125
$ws.Range("J134").Value = 3277.5
$ws.Range("K134").Value = 7389.375
$ws.Range("L134").Value = 9832.5
$ws.Range("M134").Value = -4854.375
$ws.Range("N134").Value = -14902.5

$ws = $wb.Worksheets.Item("CRP")
# Row 48 (G48 old context)
$ws.Range("H48").Value = 38890.25
$ws.Range("J48").Value = 38890.25
$ws.Range("L48").Value = 38890.25
$ws.Range("N48").Value = -39842.25

# Row 64 (G64 old context)
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("N64").Value = 0

# Row 67 (G67 old context)
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("N67").Value = 0

# Row 107 (G107 old context)
$ws.Range("H107").Value = 908.6667
$ws.Range("I107").Value = 884.3333
$ws.Range("K107").Value = 884.3333
$ws.Range("M107").Value = 1035.6667

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (G4 old context)
$ws.Range("H4").Value = 18136974
$ws.Range("I4").Value = 43062920
$ws.Range("K4").Value = 129188760
$ws.Range("M4").Value = -129188648

# Row 25 (G25 old context)
$ws.Range("H25").Value = 5167.3335
$ws.Range("J25").Value = 5251
$ws.Range("L25").Value = 15753
$ws.Range("N25").Value = -16091

# Row 30 (G30 old context)
$ws.Range("H30").Value = 5167.3335
$ws.Range("J30").Value = 5251
$ws.Range("L30").Value = 15753
$ws.Range("N30").Value = -15957

# Row 50 (G50 old context)
$ws.Range("H50").Value = 1069.2
$ws.Range("J50").Value = 781.3333
$ws.Range("L50").Value = 2343.9999
$ws.Range("N50").Value = -3305.9999

# Row 53 (G53 old context)
$ws.Range("H53").Value = 1069.2
$ws.Range("J53").Value = 781.3333
$ws.Range("L53").Value = 2343.9999
$ws.Range("N53").Value = -3305.9999

# Row 102 (G102 old context)
$ws.Range("H102").Value = 3496.5715
$ws.Range("I102").Value = 2912.8333
$ws.Range("J102").Value = 6999
$ws.Range("K102").Value = 8738.499899999999
$ws.Range("L102").Value = 20997
$ws.Range("M102").Value = -6304.499899999999
$ws.Range("N102").Value = -25865

# Row 128 (G128 old context)
$ws.Range("H128").Value = 196514.5
$ws.Range("I128").Value = 196514.5
$ws.Range("K128").Value = 589543.5
$ws.Range("M128").Value = -584563.5

# Row 133 (G133 old context)
$ws.Range("H133").Value = 6358
$ws.Range("I133").Value = 1199.5
$ws.Range("J133").Value = 11516.5
$ws.Range("K133").Value = 3598.5
$ws.Range("L133").Value = 34549.5
$ws.Range("M133").Value = 1461.5
$ws.Range("N133").Value = -44669.5

# Row 138 (G138 old context)
$ws.Range("H138").Value = 4892.3335
$ws.Range("I138").Value = 4892.3335
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 14677.0005
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -9537.000499999998

# Row 139 (G139 old context)
$ws.Range("H139").Value = 3197.7856
$ws.Range("I139").Value = 2954.5
$ws.Range("J139").Value = 3806
$ws.Range("K139").Value = 8863.5
$ws.Range("L139").Value = 11418
$ws.Range("M139").Value = -3723.5
$ws.Range("N139").Value = -21698

# Row 140 (G140 old context)
$ws.Range("H140").Value = 2955.3333
$ws.Range("I140").Value = 2816.3333
$ws.Range("J140").Value = 3094.3333
$ws.Range("K140").Value = 8448.999899999999
$ws.Range("L140").Value = 9282.999899999999
$ws.Range("M140").Value = -3268.999899999999
$ws.Range("N140").Value = -19642.9999

# Row 141 (G141 old context)
$ws.Range("H141").Value = 14257.75
$ws.Range("J141").Value = 18677.334
$ws.Range("L141").Value = 56032.00199999999
$ws.Range("N141").Value = -66392.00199999999

$ws = $wb.Worksheets.Item("GSM")
# Row 49 (G49 old context)
$ws.Range("H49").Value = 21481.25

# Row 70 (G70 old context)
$ws.Range("H70").Value = 7643.4287
$ws.Range("J70").Value = 7643.4287
$ws.Range("L70").Value = 7643.4287
$ws.Range("N70").Value = -8183.4287

# Row 73 (G73 old context)
$ws.Range("H73").Value = 7643.4287
$ws.Range("J73").Value = 7643.4287
$ws.Range("L73").Value = 7643.4287
$ws.Range("N73").Value = -9515.4287

# Row 80 (G80 old context)
$ws.Range("H80").Value = 8672.666999999999
$ws.Range("I80").Value = 9832.166999999999
$ws.Range("J80").Value = 7899.6665
$ws.Range("K80").Value = 9832.166999999999
$ws.Range("L80").Value = 7899.6665
$ws.Range("M80").Value = -8834.166999999999
$ws.Range("N80").Value = -9895.666499999999

# Row 83 (G83 old context)
$ws.Range("H83").Value = 8672.666999999999
$ws.Range("I83").Value = 9832.166999999999
$ws.Range("J83").Value = 7899.6665
$ws.Range("K83").Value = 49160.835
$ws.Range("L83").Value = 39498.3325
$ws.Range("M83").Value = -44168.835
$ws.Range("N83").Value = -49482.3325

# Row 127 (G127 old context)
$ws.Range("H127").Value = 359948.5
$ws.Range("I127").Value = 120000
$ws.Range("K127").Value = 120000
$ws.Range("M127").Value = -115040

# Row 132 (G132 old context)
$ws.Range("H132").Value = 4887.077
$ws.Range("I132").Value = 4926.696
$ws.Range("K132").Value = 14780.088
$ws.Range("M132").Value = -12250.088

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (G40 old context)
$ws.Range("H40").Value = 2935.1428
$ws.Range("I40").Value = 3309.4
$ws.Range("J40").Value = 1999.5
$ws.Range("K40").Value = 3309.4
$ws.Range("L40").Value = 1999.5
$ws.Range("M40").Value = -3173.4
$ws.Range("N40").Value = -2271.5

# Row 132 (G132 old context)
$ws.Range("H132").Value = 4653.143
$ws.Range("J132").Value = 6243.25
$ws.Range("L132").Value = 18729.75
$ws.Range("N132").Value = -23789.75

$ws = $wb.Worksheets.Item("WVR")
# Row 47 (G47 old context)
$ws.Range("H47").Value = 55822.43
$ws.Range("I47").Value = 58877
$ws.Range("J47").Value = 48186
$ws.Range("K47").Value = 58877
$ws.Range("L47").Value = 48186
$ws.Range("M47").Value = -58305
$ws.Range("N47").Value = -49330

# Row 54 (G54 old context)
$ws.Range("H54").Value = 38235
$ws.Range("J54").Value = 38235
$ws.Range("L54").Value = 38235
$ws.Range("N54").Value = -39275

# Row 132 (G132 old context)
$ws.Range("H132").Value = 1577.55
$ws.Range("I132").Value = 1634.2632
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 4902.7896
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -2372.7896
$ws.Range("N132").Value = -6560

# Row 136 (G136 old context)
$ws.Range("H136").Value = 2996.2083
$ws.Range("I136").Value = 2624.2354
$ws.Range("J136").Value = 3899.5715
$ws.Range("K136").Value = 7872.706200000001
$ws.Range("L136").Value = 11698.7145
$ws.Range("M136").Value = -5322.706200000001
$ws.Range("N136").Value = -16798.7145
